# Adjust position of multiplicity textboxes in the class diagram.
#
# PowerPoint's Shape.Left / Shape.Top are expressed in points, while the
# OOXML stores EMU (1 pt = 12700 EMU). The target positions (in EMU) are:
#   Id 54 (TextBox 53): (4480675, 3051630)
#   Id 61 (TextBox 60): (6254670, 3034551)
#   Id 65 (TextBox 64): (2590800, 2819400)
#   Id 66 (TextBox 65): (2590800, 3631317)
#   Id 70 (TextBox 69): (6619494, 3204824)
# The point values below are chosen so that round-tripping through the
# COM Single-precision Left/Top properties reproduces those exact EMU
# offsets.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$moves = @{
    54 = @{ Left = 352.80906711811025; Top = 240.28582677165355 }
    61 = @{ Left = 492.49370078740156; Top = 238.94102362204725 }
    65 = @{ Left = 204.0;              Top = 222.0 }
    66 = @{ Left = 204.0;              Top = 285.93047244094487 }
    70 = @{ Left = 521.220002;         Top = 252.3483514566929 }
}

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shape = $s.Shapes.Item($i)
    if ($moves.ContainsKey($shape.Id)) {
        $target = $moves[$shape.Id]
        $shape.Left = $target.Left
        $shape.Top = $target.Top
    }
}
